# Trade #142 (HighProbConvergence) closes as an early exit, and two new
# OPEN trades (momentum #171, MarketMaking #172) get appended to the log.
#
# Helper: write a literal text value into a cell without letting Excel's
# automatic type inference turn a date/time-looking string (e.g.
# "2026-02-18") into a real date serial number. We briefly force a Text
# number format, assign the literal value, then clear the format again so
# the cell is left with the default (unstyled) look, matching how the
# rest of the workbook's text cells are stored.
function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet — roll-up totals after the new trades / closed trade
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(6, 2).Value = 142    # Total Trades: 141 -> 142
$summary.Cells.Item(9, 2).Value = 46.48  # Win Rate %:   46.81 -> 46.48

# ---------------------------------------------------------------------
# 2) Strategy Status sheet — HighProbConvergence row (row 3)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Cells.Item(3, 4).Value = 16     # Trades:      15 -> 16
$status.Cells.Item(3, 7).Value = 68.75  # Win Rate %:  73.33 -> 68.75

# ---------------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# 3a) Existing trade #142 (row 143) closes via early_exit
$allTrades.Cells.Item(143, 7).Value = 0.09           # Exit Price
Set-TextValue $allTrades.Cells.Item(143, 8) "CLOSED" # Status
$allTrades.Cells.Item(143, 11).Value = 100.42        # Capital After
Set-TextValue $allTrades.Cells.Item(143, 12) "early_exit" # Exit Reason
$allTrades.Cells.Item(143, 13).Value = 0.16          # Duration (min)

# 3b) New trade #171 — momentum, OPEN (row 172)
$r = 172
$allTrades.Cells.Item($r, 1).Value = 171
Set-TextValue $allTrades.Cells.Item($r, 2) "2026-02-18"
Set-TextValue $allTrades.Cells.Item($r, 3) "00:38:18"
Set-TextValue $allTrades.Cells.Item($r, 4) "momentum"
Set-TextValue $allTrades.Cells.Item($r, 5) "DOWN"
$allTrades.Cells.Item($r, 6).Value = 0.09
Set-TextValue $allTrades.Cells.Item($r, 8) "OPEN"
$allTrades.Cells.Item($r, 9).Value = 0
$allTrades.Cells.Item($r, 10).Value = 0
$allTrades.Cells.Item($r, 11).Value = 99.22374292899114
$allTrades.Cells.Item($r, 13).Value = 0
$allTrades.Cells.Item($r, 14).Value = 0
$allTrades.Cells.Item($r, 15).Value = 0
$allTrades.Cells.Item($r, 16).Value = 0.9
Set-TextValue $allTrades.Cells.Item($r, 17) "Downward momentum: -1.942% over 10 samples"

# 3c) New trade #172 — MarketMaking, OPEN (row 173)
$r = 173
$allTrades.Cells.Item($r, 1).Value = 172
Set-TextValue $allTrades.Cells.Item($r, 2) "2026-02-18"
Set-TextValue $allTrades.Cells.Item($r, 3) "00:38:18"
Set-TextValue $allTrades.Cells.Item($r, 4) "MarketMaking"
Set-TextValue $allTrades.Cells.Item($r, 5) "UP"
$allTrades.Cells.Item($r, 6).Value = 0.92
Set-TextValue $allTrades.Cells.Item($r, 8) "OPEN"
$allTrades.Cells.Item($r, 9).Value = 0
$allTrades.Cells.Item($r, 10).Value = 0
$allTrades.Cells.Item($r, 11).Value = 99.20858346467945
$allTrades.Cells.Item($r, 13).Value = 0
$allTrades.Cells.Item($r, 14).Value = 0
$allTrades.Cells.Item($r, 15).Value = 0
$allTrades.Cells.Item($r, 16).Value = 0.6
Set-TextValue $allTrades.Cells.Item($r, 17) "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------
# 4) momentum sheet — same new trade #171 as above (row 45)
#    Column order differs from "All Trades":
#    L=Entry Slippage, M=Exit Slippage, N=Confidence,
#    O=Entry Reason, P=Exit Reason, Q=Duration (min)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$r = 45
$momentum.Cells.Item($r, 1).Value = 171
Set-TextValue $momentum.Cells.Item($r, 2) "2026-02-18"
Set-TextValue $momentum.Cells.Item($r, 3) "00:38:18"
Set-TextValue $momentum.Cells.Item($r, 4) "momentum"
Set-TextValue $momentum.Cells.Item($r, 5) "DOWN"
$momentum.Cells.Item($r, 6).Value = 0.09
Set-TextValue $momentum.Cells.Item($r, 8) "OPEN"
$momentum.Cells.Item($r, 9).Value = 0
$momentum.Cells.Item($r, 10).Value = 0
$momentum.Cells.Item($r, 11).Value = 99.22374292899114
$momentum.Cells.Item($r, 12).Value = 0
$momentum.Cells.Item($r, 13).Value = 0
$momentum.Cells.Item($r, 14).Value = 0.9
Set-TextValue $momentum.Cells.Item($r, 15) "Downward momentum: -1.942% over 10 samples"
$momentum.Cells.Item($r, 17).Value = 0

# ---------------------------------------------------------------------
# 5) HighProbConvergence sheet — existing trade #142 (row 17) closes
#    Column order: L=Entry Slippage, M=Exit Slippage, N=Confidence,
#    O=Entry Reason, P=Exit Reason, Q=Duration (min)
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(17, 7).Value = 0.09             # Exit Price
Set-TextValue $hpc.Cells.Item(17, 8) "CLOSED"   # Status
$hpc.Cells.Item(17, 11).Value = 100.42          # Capital After
Set-TextValue $hpc.Cells.Item(17, 16) "early_exit" # Exit Reason
$hpc.Cells.Item(17, 17).Value = 0.16            # Duration (min)

# ---------------------------------------------------------------------
# 6) MarketMaking sheet — same new trade #172 as above (row 68)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$r = 68
$mm.Cells.Item($r, 1).Value = 172
Set-TextValue $mm.Cells.Item($r, 2) "2026-02-18"
Set-TextValue $mm.Cells.Item($r, 3) "00:38:18"
Set-TextValue $mm.Cells.Item($r, 4) "MarketMaking"
Set-TextValue $mm.Cells.Item($r, 5) "UP"
$mm.Cells.Item($r, 6).Value = 0.92
Set-TextValue $mm.Cells.Item($r, 8) "OPEN"
$mm.Cells.Item($r, 9).Value = 0
$mm.Cells.Item($r, 10).Value = 0
$mm.Cells.Item($r, 11).Value = 99.20858346467945
$mm.Cells.Item($r, 12).Value = 0
$mm.Cells.Item($r, 13).Value = 0
$mm.Cells.Item($r, 14).Value = 0.6
Set-TextValue $mm.Cells.Item($r, 15) "Normal spread capture: 198 bps"
$mm.Cells.Item($r, 17).Value = 0
